$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 115
$ws.Range("F4").Value = 1829
$ws.Range("F6").Value = 3867
$ws.Range("F7").Value = 513
$ws.Range("F13").Value = 2137
$ws.Range("F15").Value = 638763
$ws.Range("F16").Value = 1570
$ws.Range("F18").Value = 1382
$ws.Range("F19").Value = 659
$ws.Range("F21").Value = 1229
$ws.Range("F22").Value = 2114
$ws.Range("F23").Value = 1084
$ws.Range("F24").Value = 2638
$ws.Range("F25").Value = 1507
$ws.Range("F26").Value = 727
$ws.Range("F27").Value = 1477
$ws.Range("F28").Value = 21
$ws.Range("F30").Value = 1061
$ws.Range("F31").Value = 225
$ws.Range("F32").Value = 1062
$ws.Range("F34").Value = 67
$ws.Range("F35").Value = 1977
$ws.Range("F36").Value = 1289
$ws.Range("F37").Value = 551
$ws.Range("F38").Value = 1184
$ws.Range("F39").Value = 1119
$ws.Range("F42").Value = 43
$ws.Range("F43").Value = 2513
$ws.Range("F44").Value = 197
$ws.Range("F45").Value = 958
$ws.Range("F46").Value = 3047
$ws.Range("F47").Value = 23

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 469
$ws.Range("F7").Value = 51
$ws.Range("F10").Value = 89
$ws.Range("F11").Value = 144109
$ws.Range("F12").Value = 144109
$ws.Range("F15").Value = 19
$ws.Range("F18").Value = 221
$ws.Range("F21").Value = 392
$ws.Range("F22").Value = 392
$ws.Range("F23").Value = 97
$ws.Range("F26").Value = 85
$ws.Range("F27").Value = 493
$ws.Range("F32").Value = 288
$ws.Range("F33").Value = 259

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 3099
$ws.Range("G5").Value = '不可售'
$ws.Range("F6").Value = 224
$ws.Range("F8").Value = 802
$ws.Range("F11").Value = 1551
$ws.Range("F12").Value = 462
$ws.Range("F13").Value = 1750

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 802
$ws.Range("F4").Value = 462
$ws.Range("F5").Value = 115
$ws.Range("F6").Value = 1829
$ws.Range("F7").Value = 1750
$ws.Range("F8").Value = 3867
$ws.Range("F9").Value = 51
$ws.Range("F10").Value = 513
$ws.Range("F14").Value = 2137
$ws.Range("F16").Value = 638764
$ws.Range("F18").Value = 89
$ws.Range("F19").Value = 1570
$ws.Range("F20").Value = 144109
$ws.Range("F22").Value = 1382
$ws.Range("F23").Value = 659
$ws.Range("F25").Value = 1229
$ws.Range("F26").Value = 2114
$ws.Range("F27").Value = 1084
$ws.Range("F28").Value = 2638
$ws.Range("F29").Value = 1507
$ws.Range("F30").Value = 727
$ws.Range("F31").Value = 19
$ws.Range("F32").Value = 1477
$ws.Range("F33").Value = 392
$ws.Range("F35").Value = 97
$ws.Range("F36").Value = 1061
$ws.Range("F37").Value = 1062
$ws.Range("F39").Value = 67
$ws.Range("F40").Value = 1977
$ws.Range("F41").Value = 1289
$ws.Range("F42").Value = 551
$ws.Range("F43").Value = 1184
$ws.Range("F44").Value = 1119
$ws.Range("F45").Value = 288
$ws.Range("F46").Value = 288
$ws.Range("F47").Value = 259
$ws.Range("F48").Value = 2513
$ws.Range("F49").Value = 197
$ws.Range("F50").Value = 958
$ws.Range("F51").Value = 3047
